$d = $word.ActiveDocument

# 1. Insert a new bold paragraph ("Play Blood Suckers 2 Free - Improved
#    Vampire-Themed Slot") right before the final (italic) paragraph.
#    (Done before the later deletion below - doing the delete first
#    confuses the insertion-point bookkeeping and corrupts the paragraph
#    that precedes the insertion point.)
$n = $d.Paragraphs.Count
$beforeLast = $d.Paragraphs($n - 1)
$insertionPoint = $beforeLast.Range.Duplicate
$insertionPoint.Collapse(0)

$newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Blood Suckers 2 Free - Improved Vampire-Themed Slot</w:t></w:r></w:p>' +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint.InsertXML($newParaXml)

# 2. Swap the final (italic) paragraph's DALLE-prompt text for the
#    meta-description blurb that used to live near the top of the doc.
$finalPara = $d.Paragraphs($d.Paragraphs.Count)
$finalPara.Range.Find.Execute(
    "Create a Feature Image Prompt for DALLE: Design a vibrant and cartoonish image featuring a happy and confident Maya warrior, wearing glasses. The warrior should appear to be holding a crossbow while standing confidently in front of a spooky background featuring luxurious coffins. Make sure to include some blood-red elements to emphasize the vampire theme of the game " + [char]34 + "Blood Suckers 2" + [char]34 + ".",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Read our review of Blood Suckers 2, an improved vampire-themed slot machine game from NetEnt. Play for free and enjoy classic Free Spins and a thrilling bonus game.",
    2)

# 3. Remove the "Meta description" paragraph that follows the H1 title
#    ("Meta description" in bold + the review blurb in a plain run).
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Delete()
